$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4.917742342958388
$ws.Range("D2").Value = 8.316818503261121
$ws.Range("E2").Value = 14.47219969740644
$ws.Range("F2").Value = 28.26786915847341
$ws.Range("G2").Value = 3.634766048845353
$ws.Range("J2").Value = 9.986270084265003
$ws.Range("M2").Value = 59.45321102715619
$ws.Range("O2").Value = 22.87115159022751
$ws.Range("C3").Value = 4.749969378258453
$ws.Range("D3").Value = 8.271155819216194
$ws.Range("E3").Value = 14.12755170890717
$ws.Range("F3").Value = 28.83485314175936
$ws.Range("G3").Value = 3.639340483171709
$ws.Range("J3").Value = 9.870200979023961
$ws.Range("M3").Value = 56.04813396013238
$ws.Range("O3").Value = 23.17593905720902
$ws.Range("C4").Value = 4.645672434795215
$ws.Range("D4").Value = 8.244978271907659
$ws.Range("E4").Value = 13.91736430526287
$ws.Range("F4").Value = 29.20196513518149
$ws.Range("G4").Value = 3.642266832652225
$ws.Range("J4").Value = 9.802418069002709
$ws.Range("M4").Value = 53.84097479413966
$ws.Range("O4").Value = 23.37744841904998
$ws.Range("C5").Value = 4.602921816618641
$ws.Range("D5").Value = 8.234785319342253
$ws.Range("E5").Value = 13.83219487923201
$ws.Range("F5").Value = 29.35626824619308
$ws.Range("G5").Value = 3.643489124441372
$ws.Range("J5").Value = 9.775693631376454
$ws.Range("M5").Value = 52.9125288738847
$ws.Range("O5").Value = 23.46310542460593
$ws.Range("C6").Value = 4.595810354940393
$ws.Range("D6").Value = 8.233121650500355
$ws.Range("E6").Value = 13.81808548657684
$ws.Range("F6").Value = 29.38217218941338
$ws.Range("G6").Value = 3.643693889553909
$ws.Range("J6").Value = 9.771310875515486
$ws.Range("M6").Value = 52.7566184140152
$ws.Range("O6").Value = 23.47754046780102
$ws.Range("C7").Value = 4.645096790382016
$ws.Range("D7").Value = 8.244838875492615
$ws.Range("E7").Value = 13.91621355486155
$ws.Range("F7").Value = 29.20402717193308
$ws.Range("G7").Value = 3.642283196061312
$ws.Range("J7").Value = 9.80205399192698
$ws.Range("M7").Value = 53.82857049805555
$ws.Range("O7").Value = 23.37858937842231
$ws.Range("C8").Value = 4.860203322702902
$ws.Range("D8").Value = 8.300692375637214
$ws.Range("E8").Value = 14.35314217652524
$ws.Range("F8").Value = 28.45938092798213
$ws.Range("G8").Value = 3.636319023726463
$ws.Range("J8").Value = 9.945539719150345
$ws.Range("M8").Value = 58.30343977742112
$ws.Range("O8").Value = 22.97321525337693
$ws.Range("C9").Value = 5.268880682295394
$ws.Range("D9").Value = 8.424655089146906
$ws.Range("E9").Value = 15.21605777287157
$ws.Range("F9").Value = 27.15326519261207
$ws.Range("G9").Value = 3.625547090640773
$ws.Range("J9").Value = 10.25355341165983
$ws.Range("M9").Value = 66.14702533219976
$ws.Range("O9").Value = 22.29572318721957
$ws.Range("C10").Value = 5.557578396818967
$ws.Range("D10").Value = 8.524087795835781
$ws.Range("E10").Value = 15.84719759033394
$ws.Range("F10").Value = 26.29306517421236
$ws.Range("G10").Value = 3.618182283603474
$ws.Range("J10").Value = 10.49468919468748
$ws.Range("M10").Value = 71.33568839428135
$ws.Range("O10").Value = 21.87456701208706
$ws.Range("C11").Value = 5.685775116941544
$ws.Range("D11").Value = 8.571039004193649
$ws.Range("E11").Value = 16.13245860968845
$ws.Range("F11").Value = 25.92471066894415
$ws.Range("G11").Value = 3.614948005816924
$ws.Range("J11").Value = 10.60730229053348
$ws.Range("M11").Value = 73.57118553681846
$ws.Range("O11").Value = 21.70077385544278
$ws.Range("C12").Value = 5.733824541540144
$ws.Range("D12").Value = 8.589056811257262
$ws.Range("E12").Value = 16.24012103976469
$ws.Range("F12").Value = 25.78865864299384
$ws.Range("G12").Value = 3.613739695468108
$ws.Range("J12").Value = 10.65033965540217
$ws.Range("M12").Value = 74.39977112940458
$ws.Range("O12").Value = 21.63762527776229
$ws.Range("C13").Value = 5.723498957713913
$ws.Range("D13").Value = 8.585165890703305
$ws.Range("E13").Value = 16.21695141337405
$ws.Range("F13").Value = 25.81780476494914
$ws.Range("G13").Value = 3.613999199460633
$ws.Range("J13").Value = 10.64105371804385
$ws.Range("M13").Value = 74.22211829261951
$ws.Range("O13").Value = 21.65110529001717
$ws.Range("C14").Value = 5.689738350883925
$ws.Range("D14").Value = 8.572516612367867
$ws.Range("E14").Value = 16.14132378007207
$ws.Range("F14").Value = 25.91344790254555
$ws.Range("G14").Value = 3.614848269319877
$ws.Range("J14").Value = 10.61083527021187
$ws.Range("M14").Value = 73.63971336170822
$ws.Range("O14").Value = 21.69552459366258
$ws.Range("C15").Value = 5.668993171659044
$ws.Range("D15").Value = 8.564799350914512
$ws.Range("E15").Value = 16.0949501802551
$ws.Range("F15").Value = 25.97248378793777
$ws.Range("G15").Value = 3.615370482848892
$ws.Range("J15").Value = 10.59237602725283
$ws.Range("M15").Value = 73.28063670869524
$ws.Range("O15").Value = 21.72308270883185
$ws.Range("C16").Value = 5.549133390851542
$ws.Range("D16").Value = 8.521053262612693
$ws.Range("E16").Value = 15.8285096974629
$ws.Range("F16").Value = 26.31761112374593
$ws.Range("G16").Value = 3.618395965203572
$ws.Range("J16").Value = 10.48738612834723
$ws.Range("M16").Value = 71.18708081678244
$ws.Range("O16").Value = 21.88629165638683
$ws.Range("C17").Value = 5.47476679638054
$ws.Range("D17").Value = 8.494650459526371
$ws.Range("E17").Value = 15.66451423557874
$ws.Range("F17").Value = 26.53530634597909
$ws.Range("G17").Value = 3.620281544102542
$ws.Range("J17").Value = 10.42370691482705
$ws.Range("M17").Value = 69.87076864998119
$ws.Range("O17").Value = 21.99104494511542
$ws.Range("C18").Value = 5.431700528591175
$ws.Range("D18").Value = 8.479626573225167
$ws.Range("E18").Value = 15.57001956148034
$ws.Range("F18").Value = 26.66267057612285
$ws.Range("G18").Value = 3.62137701354099
$ws.Range("J18").Value = 10.38735666405402
$ws.Range("M18").Value = 69.10190673537471
$ws.Range("O18").Value = 22.05296444057391
$ws.Range("C19").Value = 5.417070251505733
$ws.Range("D19").Value = 8.474567885982481
$ws.Range("E19").Value = 15.53799929914085
$ws.Range("F19").Value = 26.70615959402782
$ws.Range("G19").Value = 3.62174980588394
$ws.Range("J19").Value = 10.37509737140078
$ws.Range("M19").Value = 68.83956487744814
$ws.Range("O19").Value = 22.07421272990515
$ws.Range("C20").Value = 5.482713872680358
$ws.Range("D20").Value = 8.497444348056646
$ws.Range("E20").Value = 15.68199001345154
$ws.Range("F20").Value = 26.51190869408706
$ws.Range("G20").Value = 3.620079691074257
$ws.Range("J20").Value = 10.43045726830749
$ws.Range("M20").Value = 70.01210887424531
$ws.Range("O20").Value = 21.97972038873253
$ws.Range("C21").Value = 5.699668452368846
$ws.Range("D21").Value = 8.5762256116787
$ws.Range("E21").Value = 16.16354791860728
$ws.Range("F21").Value = 25.88526079829321
$ws.Range("G21").Value = 3.61459843269529
$ws.Range("J21").Value = 10.61970069818087
$ws.Range("M21").Value = 73.811266750857
$ws.Range("O21").Value = 21.68240443171767
$ws.Range("C22").Value = 5.838553082351649
$ws.Range("D22").Value = 8.629099255633776
$ws.Range("E22").Value = 16.47613960308367
$ws.Range("F22").Value = 25.49581025068255
$ws.Range("G22").Value = 3.611111823182579
$ws.Range("J22").Value = 10.74566230826805
$ws.Range("M22").Value = 76.18963869571334
$ws.Range("O22").Value = 21.50366624567863
$ws.Range("C23").Value = 5.764707397592246
$ws.Range("D23").Value = 8.600755715557003
$ws.Range("E23").Value = 16.30952764660704
$ws.Range("F23").Value = 25.70177969055787
$ws.Range("G23").Value = 3.612964016768765
$ws.Range("J23").Value = 10.67823424673887
$ws.Range("M23").Value = 74.9298188137013
$ws.Range("O23").Value = 21.59760115395501
$ws.Range("C24").Value = 5.47912196755574
$ws.Range("D24").Value = 8.49618074623705
$ws.Range("E24").Value = 15.67408986581282
$ws.Range("F24").Value = 26.52247990402518
$ws.Range("G24").Value = 3.620170913223572
$ws.Range("J24").Value = 10.42740462233064
$ws.Range("M24").Value = 69.94824663686799
$ws.Range("O24").Value = 21.98483494762278
$ws.Range("C25").Value = 5.160119058634513
$ws.Range("D25").Value = 8.389619303891406
$ws.Range("E25").Value = 14.98266726669127
$ws.Range("F25").Value = 27.48957853842506
$ws.Range("G25").Value = 3.628363690831243
$ws.Range("J25").Value = 10.16751768712346
$ws.Range("M25").Value = 64.12568954370343
$ws.Range("O25").Value = 22.46593039346013
